$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "tree" model parameters: min_samples_split range 1->2
$ws.Range("B2").Value = "{'max_depth': [2, 100], 'min_samples_split': [2, 100], 'min_samples_leaf': [2, 100]}"

# Header: rename "max_lags" column to "lags"
$ws.Range("C1").Value = "lags"

# Replace the max_lags numeric column (10) with the new lags list for every model row
$ws.Range("C2").Value = "[0,5,10,15]"
$ws.Range("C3").Value = "[0,5,10,15]"
$ws.Range("C4").Value = "[0,5,10,15]"
$ws.Range("C5").Value = "[0,5,10,15]"
$ws.Range("C6").Value = "[0,5,10,15]"
$ws.Range("C7").Value = "[0,5,10,15]"

# Move the active selection to B5
$ws.Range("B5").Select()
